$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6860853433609009
$ws.Range("B1").Value = 1.041175246238708
$ws.Range("C1").Value = 2.278500318527222
$ws.Range("D1").Value = 3.766931533813477
$ws.Range("E1").Value = 1.644686818122864
